# Add a new "2021" column (R) to the table: copy the formatting from the
# existing "2020" column (Q) for each of the header/data rows, then fill
# in the new year header and data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NewYearCell($srcAddr, $dstAddr, $value) {
    $ws.Range($srcAddr).Copy() | Out-Null
    $ws.Range($dstAddr).PasteSpecial(-4122) | Out-Null  # xlPasteFormats
    $ws.Range($dstAddr).Value = $value
}

Set-NewYearCell "Q4" "R4" 2021
Set-NewYearCell "Q5" "R5" 47.8
Set-NewYearCell "Q6" "R6" 20.7
Set-NewYearCell "Q7" "R7" 9.8
Set-NewYearCell "Q8" "R8" 17.3

$excel.CutCopyMode = $false

# Move the active selection to P10, matching the saved view state
$ws.Range("P10").Select() | Out-Null
